$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto-price snapshot refresh (GitHub Actions bot).
# A handful of the new "Price" values parse as plain numbers; force
# those cells to Text format first so Excel keeps the literal digits
# (matching the scraped formatting, incl. trailing zeros) instead of
# silently converting them to numeric cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.843.15'
$ws.Range("D3").Value = '1.637.76'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '215.65'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").Value = '0.5064'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D8").Value = '0.2576'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").Value = '0.06427'
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").Value = '0.07769'
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").Value = '4.281'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '1.863.06'
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").Value = '1.633.19'
$ws.Range("E14").Value = '  +0.03%  '
$ws.Range("D15").Value = '0.5637'
$ws.Range("E15").Value = '  +3.82%  '
$ws.Range("D16").Value = '0.0₅7592'
$ws.Range("E16").Value = '  -1.81%  '
$ws.Range("D17").Value = '63.13'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = '25.865.12'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").Value = '194.93'
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("E21").Value = '  -2.47%  '
$ws.Range("D22").Value = '9.875'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '6.096'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -5.00%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1271'
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '139.90'
$ws.Range("E27").Value = '  -2.11%  '
$ws.Range("D28").Value = '6.783'
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").Value = '15.51'
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("D30").Value = '1.244'
$ws.Range("E30").Value = '  +0.69%  '
$ws.Range("D31").Value = '0.04868'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '3.296'
$ws.Range("E32").Value = '  +1.82%  '
$ws.Range("D33").Value = '3.219'
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("D34").Value = '1.558'
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("D35").Value = '2.376'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = '0.9035'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").Value = '2.579'
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("D38").Value = '1.132.71'
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("D39").Value = '0.5509'
$ws.Range("E39").Value = '  +0.37%  '
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("D41").Value = '0.9961'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = '5.521'
$ws.Range("E42").Value = '  -1.02%  '
$ws.Range("D43").Value = '0.8001'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = '97.78'
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").Value = '1.773.22'
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").Value = '  -7.94%  '
$ws.Range("D47").Value = '55.36'
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").Value = '0.4386'
$ws.Range("E48").Value = '  -2.09%  '
$ws.Range("D49").Value = '0.05050'
$ws.Range("E49").Value = '  -2.37%  '
$ws.Range("D50").Value = '7.672'
$ws.Range("E50").Value = '  +2.47%  '
$ws.Range("D51").Value = '1.003'
$ws.Range("E51").Value = '  -0.18%  '
